# Applies the "Updated cryptos list" data refresh to Sheet1.
#
# Columns: A=rank(idx) B=Coin C=Link D=Price E=Volume(1h)
# Many Price values look like plain numbers (e.g. "239.33") but must stay
# stored as TEXT (matching the source data, which used inline/shared
# strings, not numeric cells). Assigning such a string straight to
# Range.Value makes Excel auto-convert it to a number, so Set-CellText
# below forces text formatting for the duration of the write and then
# clears the formatting override again so the cell's style index is left
# exactly as it was before (no stray numFmt/style differences).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText([string]$cellRef, [string]$text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

$ws.Range("D2").Value = "95.996.81"
$ws.Range("E2").Value = "  +2.46%  "

$ws.Range("D3").Value = "3.562.56"
$ws.Range("E3").Value = "  +7.12%  "

$ws.Range("E4").Value = "  +0.03%  "

Set-CellText "D5" "239.33"
$ws.Range("E5").Value = "  +3.87%  "

Set-CellText "D6" "637.09"
$ws.Range("E6").Value = "  +3.23%  "

$ws.Range("E7").Value = "  +7.28%  "

Set-CellText "D8" "0.401"
$ws.Range("E8").Value = "  +3.55%  "

$ws.Range("E9").Value = "  -0.04%  "

$ws.Range("E10").Value = "  +10.28%  "

$ws.Range("D11").Value = "3.560.53"
$ws.Range("E11").Value = "  +7.14%  "

Set-CellText "D12" "43.32"
$ws.Range("E12").Value = "  +3.44%  "

Set-CellText "D13" "0.201"
$ws.Range("E13").Value = "  +3.89%  "

Set-CellText "D14" "6.44"
$ws.Range("E14").Value = "  +8.40%  "

$ws.Range("D15").Value = "4.240.85"
$ws.Range("E15").Value = "  +7.58%  "

$ws.Range("D16").Value = "95.937.88"
$ws.Range("E16").Value = "  +2.54%  "

$ws.Range("E17").Value = "  +4.41%  "

$ws.Range("D18").Value = "3.557.25"
$ws.Range("E18").Value = "  +7.13%  "

$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-CellText "D19" "13.11"
$ws.Range("E19").Value = "  +20.58%  "

$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-CellText "D20" "7.75"
$ws.Range("E20").Value = "  -3.87%  "

Set-CellText "D21" "18.01"
$ws.Range("E21").Value = "  +5.17%  "

Set-CellText "D22" "0.504"
$ws.Range("E22").Value = "  +13.08%  "

Set-CellText "D23" "516.60"
$ws.Range("E23").Value = "  +4.66%  "

Set-CellText "D24" "3.41"
$ws.Range("E24").Value = "  -1.25%  "

$ws.Range("E26").Value = "  +8.79%  "

Set-CellText "D27" "97.01"
$ws.Range("E27").Value = "  +8.30%  "

Set-CellText "D28" "12.35"
$ws.Range("E28").Value = "  +5.50%  "

Set-CellText "D29" "3.10"
$ws.Range("E29").Value = "  +18.92%  "

Set-CellText "D30" "0.146"
$ws.Range("E30").Value = "  +4.27%  "

Set-CellText "D31" "11.54"
$ws.Range("E31").Value = "  +4.68%  "

$ws.Range("E32").Value = "  -0.03%  "

Set-CellText "D33" "0.183"
$ws.Range("E33").Value = "  +5.44%  "

Set-CellText "D34" "0.999"
$ws.Range("E34").Value = "  +0.47%  "

Set-CellText "D35" "30.23"
$ws.Range("E35").Value = "  +7.05%  "

Set-CellText "D36" "0.564"
$ws.Range("E36").Value = "  +6.76%  "

Set-CellText "D37" "577.84"
$ws.Range("E37").Value = "  +9.31%  "

$ws.Range("E38").Value = "  +6.65%  "

$ws.Range("E39").Value = "  +9.45%  "

$ws.Range("E40").Value = "  -0.01%  "

$ws.Range("E41").Value = "  +2.82%  "

Set-CellText "D42" "0.922"
$ws.Range("E42").Value = "  +7.37%  "

$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-CellText "D43" "0.0433"
$ws.Range("E43").Value = "  +4.45%  "

$ws.Range("B44").Value = "ImmutableX"
$ws.Range("C44").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-CellText "D44" "1.75"
$ws.Range("E44").Value = "  +4.58%  "

Set-CellText "D45" "23.84"
$ws.Range("E45").Value = "  -0.74%  "

$ws.Range("E46").Value = "  +4.68%  "

$ws.Range("E47").Value = "  -2.87%  "

$ws.Range("E48").Value = "  +3.59%  "

Set-CellText "D49" "53.91"
$ws.Range("E49").Value = "  +3.21%  "

$ws.Range("E50").Value = "  +2.90%  "

$ws.Range("E51").Value = "  +2.70%  "
